# Update crypto price/volume data per commit:
# "Updated cryptos list on Tue May 23 15:08:21 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.399.94'
$ws.Range("D3").Value = '1.861.82'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.06'
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4627'
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3716'
$ws.Range("E8").Value = '  +0.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07317'
$ws.Range("E9").Value = '  -0.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8889'
$ws.Range("E10").Value = '  +1.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.03'
$ws.Range("E11").Value = '  +0.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07808'
$ws.Range("E12").Value = '  -2.06%  '
$ws.Range("D13").Value = '1.783.35'
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.564'
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.92'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008954'
$ws.Range("E18").Value = '  +0.89%  '
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("E20").Value = '  +0.70%  '
$ws.Range("D21").Value = '27.414.98'
$ws.Range("E21").Value = '  +1.98%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("D24").Value = '2.071.80'
$ws.Range("E24").Value = '  +4.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.931'
$ws.Range("E25").Value = '  +5.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.33'
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.41'
$ws.Range("E27").Value = '  -1.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.057'
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.34'
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08854'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("E32").Value = '  +5.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7698'
$ws.Range("E33").Value = '  +4.90%  '
$ws.Range("E34").Value = '  +3.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.515'
$ws.Range("E35").Value = '  +1.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.734'
$ws.Range("E36").Value = '  +11.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.082'
$ws.Range("E37").Value = '  +0.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01959'
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05248'
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.963'
$ws.Range("E40").Value = '  +0.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.074'
$ws.Range("E41").Value = '  -1.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5145'
$ws.Range("E42").Value = '  -0.95%  '
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.428'
$ws.Range("E44").Value = '  +2.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4812'
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.36'
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '103.28'
$ws.Range("E48").Value = '  +0.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.654'
$ws.Range("E49").Value = '  +1.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06223'
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '65.56'
$ws.Range("E51").Value = '  +1.26%  '
